$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3 each describe one species observation record. In this
# edit the two records swap places (row 2's record moves to row 3, and
# row 3's record moves to row 2). Read every value first (Value2 - plain
# data, no formatting) so the swap doesn't clobber itself, then write the
# values back out in swapped order.

$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$e2 = $ws.Range("E2").Value2
$f2 = $ws.Range("F2").Value2
$g2 = $ws.Range("G2").Value2
$h2 = $ws.Range("H2").Value2
$q2 = $ws.Range("Q2").Value2
$r2 = $ws.Range("R2").Value2

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2
$q3 = $ws.Range("Q3").Value2
$r3 = $ws.Range("R3").Value2

$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("E2").Value = $e3
$ws.Range("F2").Value = $f3
$ws.Range("G2").Value = $g3
$ws.Range("H2").Value = $h3
$ws.Range("Q2").Value = $q3
$ws.Range("R2").Value = $r3

$ws.Range("A3").Value = $a2
$ws.Range("B3").Value = $b2
$ws.Range("E3").Value = $e2
$ws.Range("F3").Value = $f2
$ws.Range("G3").Value = $g2
$ws.Range("H3").Value = $h2
$ws.Range("Q3").Value = $q2
$ws.Range("R3").Value = $r2

# Column I ("Antal") held "30" as text for the old row 3 record and was
# blank (present-but-empty text cell) for the old row 2 record. Moving the
# records swaps those too. Writing a bare numeric-looking string would get
# auto-converted to a real number, so force text via a leading apostrophe,
# then drop back to the Normal style so no stray number-format style is
# left behind on the cell.
$ws.Range("I2").Value = "'30"
$ws.Range("I2").Style = "Normal"

$ws.Range("I3").Value = "'"
$ws.Range("I3").Style = "Normal"

# Column AF ("Bestämningsmetod") was a present-but-empty text cell on row 3
# and entirely absent on row 2. After the swap it is row 2 that carries the
# present-but-empty text cell, and row 3 goes back to having no cell there.
$ws.Range("AF2").Value = "'"
$ws.Range("AF2").Style = "Normal"

$ws.Range("AF3").ClearContents()
